$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 457.5
$ws.Range("I58").Value = 225
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 675
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -525
$ws.Range("N58").Value = -3300
$ws.Range("H69").Value = 4022.2654
$ws.Range("I69").Value = 4192.857
$ws.Range("J69").Value = 3954.0286
$ws.Range("K69").Value = 12578.571
$ws.Range("L69").Value = 11862.0858
$ws.Range("M69").Value = -11704.571
$ws.Range("N69").Value = -13610.0858
$ws.Range("H72").Value = 4022.2654
$ws.Range("I72").Value = 4192.857
$ws.Range("J72").Value = 3954.0286
$ws.Range("K72").Value = 37735.713
$ws.Range("L72").Value = 35586.2574
$ws.Range("M72").Value = -33367.713
$ws.Range("N72").Value = -44322.2574
$ws.Range("H74").Value = 2996.4333
$ws.Range("I74").Value = 2613.9524
$ws.Range("J74").Value = 3888.889
$ws.Range("K74").Value = 2613.9524
$ws.Range("L74").Value = 3888.889
$ws.Range("M74").Value = -1677.9524
$ws.Range("N74").Value = -5760.889
$ws.Range("H77").Value = 2996.4333
$ws.Range("I77").Value = 2613.9524
$ws.Range("J77").Value = 3888.889
$ws.Range("K77").Value = 13069.762
$ws.Range("L77").Value = 19444.445
$ws.Range("M77").Value = -8389.762000000001
$ws.Range("N77").Value = -28804.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 83.3
$ws.Range("I4").Value = 79.125
$ws.Range("K4").Value = 79.125
$ws.Range("M4").Value = 36.875
$ws.Range("H5").Value = 356.82352
$ws.Range("I5").Value = 481.58334
$ws.Range("J5").Value = 57.4
$ws.Range("K5").Value = 481.58334
$ws.Range("L5").Value = 57.4
$ws.Range("M5").Value = -369.58334
$ws.Range("N5").Value = -281.4
$ws.Range("H63").Value = 1700.1
$ws.Range("I63").Value = 1744.6666
$ws.Range("J63").Value = 1299
$ws.Range("K63").Value = 1744.6666
$ws.Range("L63").Value = 1299
$ws.Range("M63").Value = -1058.6666
$ws.Range("N63").Value = -2671
$ws.Range("H66").Value = 1700.1
$ws.Range("I66").Value = 1744.6666
$ws.Range("J66").Value = 1299
$ws.Range("K66").Value = 8723.333000000001
$ws.Range("L66").Value = 6495
$ws.Range("M66").Value = -5291.333000000001
$ws.Range("N66").Value = -13359
$ws.Range("H97").Value = 3217.8333
$ws.Range("I97").Value = 1451.75
$ws.Range("J97").Value = 6750
$ws.Range("K97").Value = 1451.75
$ws.Range("L97").Value = 6750
$ws.Range("M97").Value = -955.75
$ws.Range("N97").Value = -7742
$ws.Range("H102").Value = 1556.7693
$ws.Range("I102").Value = 1250
$ws.Range("J102").Value = 1914.6666
$ws.Range("K102").Value = 1250
$ws.Range("L102").Value = 1914.6666
$ws.Range("M102").Value = 372
$ws.Range("N102").Value = -5158.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 356.82352
$ws.Range("I4").Value = 481.58334
$ws.Range("J4").Value = 57.4
$ws.Range("K4").Value = 481.58334
$ws.Range("L4").Value = 57.4
$ws.Range("M4").Value = -366.58334
$ws.Range("N4").Value = -287.4
$ws.Range("H22").Value = 375.5
$ws.Range("I22").Value = 400.42856
$ws.Range("J22").Value = 201
$ws.Range("K22").Value = 400.42856
$ws.Range("L22").Value = 201
$ws.Range("M22").Value = -227.42856
$ws.Range("N22").Value = -547
$ws.Range("H94").Value = 1352.7858
$ws.Range("I94").Value = 1264.5385
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1264.5385
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -813.5385000000001
$ws.Range("N94").Value = -3402
$ws.Range("H99").Value = 1786.3636
$ws.Range("I99").Value = 1183.3334
$ws.Range("J99").Value = 2012.5
$ws.Range("K99").Value = 1183.3334
$ws.Range("L99").Value = 2012.5
$ws.Range("M99").Value = 314.6666
$ws.Range("N99").Value = -5008.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 61.03125
$ws.Range("I7").Value = 32.52941
$ws.Range("J7").Value = 93.333336
$ws.Range("K7").Value = 32.52941
$ws.Range("L7").Value = 93.333336
$ws.Range("M7").Value = 80.47059
$ws.Range("N7").Value = -319.333336
$ws.Range("H22").Value = 9429.546
$ws.Range("I22").Value = 14647.714
$ws.Range("J22").Value = 297.75
$ws.Range("K22").Value = 14647.714
$ws.Range("L22").Value = 297.75
$ws.Range("M22").Value = -14297.714
$ws.Range("N22").Value = -997.75
$ws.Range("H62").Value = 4760
$ws.Range("I62").Value = 5100
$ws.Range("J62").Value = 2975
$ws.Range("K62").Value = 5100
$ws.Range("L62").Value = 2975
$ws.Range("M62").Value = -4476
$ws.Range("N62").Value = -4223
$ws.Range("H65").Value = 4760
$ws.Range("I65").Value = 5100
$ws.Range("J65").Value = 2975
$ws.Range("K65").Value = 25500
$ws.Range("L65").Value = 14875
$ws.Range("M65").Value = -22380
$ws.Range("N65").Value = -21115
$ws.Range("H97").Value = 34000
$ws.Range("J97").Value = 34000
$ws.Range("L97").Value = 34000
$ws.Range("N97").Value = -35982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2213.4
$ws.Range("I17").Value = 116.666664
$ws.Range("J17").Value = 3611.2222
$ws.Range("K17").Value = 349.999992
$ws.Range("L17").Value = 10833.6666
$ws.Range("M17").Value = -180.999992
$ws.Range("N17").Value = -11171.6666
$ws.Range("H34").Value = 1375.7858
$ws.Range("I34").Value = 325
$ws.Range("J34").Value = 1550.9166
$ws.Range("K34").Value = 975
$ws.Range("L34").Value = 4652.7498
$ws.Range("M34").Value = -891
$ws.Range("N34").Value = -4820.7498
$ws.Range("H39").Value = 2935
$ws.Range("J39").Value = 3282.8572
$ws.Range("L39").Value = 9848.571599999999
$ws.Range("N39").Value = -10436.5716
$ws.Range("H55").Value = 2842.2
$ws.Range("J55").Value = 3177.75
$ws.Range("L55").Value = 9533.25
$ws.Range("N55").Value = -9887.25
$ws.Range("H63").Value = 1975
$ws.Range("I63").Value = 300
$ws.Range("J63").Value = 3650
$ws.Range("K63").Value = 900
$ws.Range("L63").Value = 10950
$ws.Range("M63").Value = -151
$ws.Range("N63").Value = -12448
$ws.Range("H64").Value = 1300
$ws.Range("I64").Value = 400
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 1200
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -930
$ws.Range("N64").Value = -12540
$ws.Range("H66").Value = 1975
$ws.Range("I66").Value = 300
$ws.Range("J66").Value = 3650
$ws.Range("K66").Value = 2700
$ws.Range("L66").Value = 32850
$ws.Range("M66").Value = 1044
$ws.Range("N66").Value = -40338
$ws.Range("H67").Value = 1300
$ws.Range("I67").Value = 400
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 1200
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -264
$ws.Range("N67").Value = -13872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28.5
$ws.Range("I2").Value = 22.4
$ws.Range("J2").Value = 38.666668
$ws.Range("K2").Value = 22.4
$ws.Range("L2").Value = 38.666668
$ws.Range("M2").Value = 90.59999999999999
$ws.Range("N2").Value = -264.666668
$ws.Range("H32").Value = 36145
$ws.Range("J32").Value = 36145
$ws.Range("L32").Value = 36145
$ws.Range("N32").Value = -36737
$ws.Range("H45").Value = 29666.666
$ws.Range("I45").Value = 20000
$ws.Range("J45").Value = 34500
$ws.Range("K45").Value = 20000
$ws.Range("L45").Value = 34500
$ws.Range("M45").Value = -19441
$ws.Range("N45").Value = -35618
$ws.Range("H80").Value = 2420.4333
$ws.Range("I80").Value = 1854.3158
$ws.Range("J80").Value = 3398.2727
$ws.Range("K80").Value = 1854.3158
$ws.Range("L80").Value = 3398.2727
$ws.Range("M80").Value = -856.3158000000001
$ws.Range("N80").Value = -5394.2727
$ws.Range("H83").Value = 2420.4333
$ws.Range("I83").Value = 1854.3158
$ws.Range("J83").Value = 3398.2727
$ws.Range("K83").Value = 9271.579
$ws.Range("L83").Value = 16991.3635
$ws.Range("M83").Value = -4279.579
$ws.Range("N83").Value = -26975.3635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1696.8125
$ws.Range("I46").Value = 2338.1428
$ws.Range("J46").Value = 1198
$ws.Range("K46").Value = 2338.1428
$ws.Range("L46").Value = 1198
$ws.Range("M46").Value = -2150.1428
$ws.Range("N46").Value = -1574
$ws.Range("H68").Value = 1669.862
$ws.Range("J68").Value = 1928.8889
$ws.Range("L68").Value = 1928.8889
$ws.Range("N68").Value = -3426.8889
$ws.Range("H71").Value = 1669.862
$ws.Range("J71").Value = 1928.8889
$ws.Range("L71").Value = 9644.4445
$ws.Range("N71").Value = -17132.4445
$ws.Range("H140").Value = 29950
$ws.Range("J140").Value = 29950
$ws.Range("L140").Value = 29950
$ws.Range("N140").Value = -40310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 7000
$ws.Range("J44").Value = 7000
$ws.Range("L44").Value = 7000
$ws.Range("N44").Value = -8108
$ws.Range("H62").Value = 2776.7368
$ws.Range("I62").Value = 2773.4546
$ws.Range("J62").Value = 2781.25
$ws.Range("K62").Value = 2773.4546
$ws.Range("L62").Value = 2781.25
$ws.Range("M62").Value = -2149.4546
$ws.Range("N62").Value = -4029.25
$ws.Range("H65").Value = 2776.7368
$ws.Range("I65").Value = 2773.4546
$ws.Range("J65").Value = 2781.25
$ws.Range("K65").Value = 13867.273
$ws.Range("L65").Value = 13906.25
$ws.Range("M65").Value = -10747.273
$ws.Range("N65").Value = -20146.25
$ws.Range("H96").Value = 2629
$ws.Range("I96").Value = 1933
$ws.Range("J96").Value = 3151
$ws.Range("K96").Value = 1933
$ws.Range("L96").Value = 3151
$ws.Range("M96").Value = -560
$ws.Range("N96").Value = -5897
